# "Generate Report for Handback"
# For each localized-language sheet (zh-cn, de-de), mark the handed-off
# files as handed back: update the Status text, record the Latest Target
# File / Latest Handback File (with hyperlinks) and stamp the Latest
# Handback DateTime.

$wb = $excel.ActiveWorkbook

$langs = @(
  @{ Sheet = "zh-cn"; Code = "zh-cn"; XlfFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandoffCommit = "56e7b7adba75c1c7794280bbf043ced8022b9ce2"; HandbackDateTime = "2016-01-25 07:48:31" },
  @{ Sheet = "de-de"; Code = "de-de"; XlfFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandoffCommit = "9e2fe7e7cb1c53161534cc5535861171a5a8a251"; HandbackDateTime = "2016-01-25 07:48:51" }
)

$srcCommit = "f11acd9490e3900c36ac9a6ea4467c7386165796"
$statusText = "Handed back: in sync with en-US"

# The Overview sheet mirrors the same "Ready for handoff" status text for
# every language column; it needs to be brought in sync as well.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    $overview.Cells.Item($row, 2).Value = $statusText   # B: zh-cn
    $overview.Cells.Item($row, 3).Value = $statusText   # C: de-de
}

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    $targetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/a.md"
    $handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($lang.HandoffCommit)/ol-handback/OpenLocalizationTestOrg/oltest.$($lang.Code)/yuwzho/$($lang.XlfFile)"

    # Row 2 (a.md) and Row 3 (b.md) were both handed off together and are
    # reported back in this pass.
    foreach ($row in 2, 3) {
        $ws.Cells.Item($row, 2).Value = $statusText          # B: Status

        # E: Latest Target File (the localized a.md), F: Latest Handback File (xlf)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $targetUrl, "", "", "a.md")
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $handbackUrl, "", "", $lang.XlfFile)

        $ws.Cells.Item($row, 7).Value = $lang.HandbackDateTime   # G: Latest Handback DateTime
    }
}
